# Generate Report for Handback
# Update the handoff/handback generation timestamps for the
# "b5fa3b00-7fc8-45f1-91f6-e18accf757cd" row across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b5fa3b00... (row 3)
$overview.Range("G3").Value = "2016-10-21 00:33:47"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) for row 3
$zhcn.Range("H3").Value = "2016-10-21 00:33:35"
$zhcn.Range("K3").Value = "2016-10-21 00:34:35"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) for row 3
$dede.Range("H3").Value = "2016-10-21 00:33:47"
$dede.Range("K3").Value = "2016-10-21 00:34:53"
